# Update Leave Card 12/22/2023 10:59 AM
# Applies the recorded edits to BATHAN, ELVIRA.xlsx:
#  - 2018 LEAVE CREDITS (sheet1 / Table13): new monthly rows, new leave entries, table grows by one row
#  - 2017 LEAVE BALANCE (sheet2 / Table1): new leave entries in rows 36-39
#  - Active tab moves from "2017 LEAVE BALANCE" back to "2018 LEAVE CREDITS"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 2018 LEAVE CREDITS
$ws2 = $wb.Worksheets.Item(2)   # 2017 LEAVE BALANCE

# ---------------------------------------------------------------------------
# 2017 LEAVE BALANCE (sheet2) - rows 36-39 get new leave entries
# Order of first-use of brand-new text matches the target shared-string order:
#   94 VL(1-0-0), 98 11/22,23,24,28,29/2023, 99 VL(6-0-0), 100 12/21,22,26-29/2023
# ---------------------------------------------------------------------------

$ws2.Range("A36").Value = 45200
$ws2.Range("B36").Value = "VL(1-0-0)"
$ws2.Range("D36").Value = 1
$ws2.Range("K36").Value = 45233

$ws2.Range("B37").Value = "VL(1-0-0)"
$ws2.Range("D37").Value = 1
$ws2.Range("K36").Copy()
$ws2.Range("K37").PasteSpecial(-4122)
$ws2.Range("K37").Value = 45230

$ws2.Range("A38").Value = 45231
$ws2.Range("B38").Value = "VL(5-0-0)"
$ws2.Range("D38").Value = 5
$ws2.Range("K38").Value = "11/22,23,24,28,29/2023"

$ws2.Range("B39").Value = "VL(6-0-0)"
$ws2.Range("D39").Value = 6
$ws2.Range("K39").Value = "12/21,22,26-29/2023"

# ---------------------------------------------------------------------------
# 2018 LEAVE CREDITS (sheet1) - rows 89-103
# ---------------------------------------------------------------------------

$ws1.Range("C89").Value = 1.25

$ws1.Range("B90").Value = "SP(2-0-0)"
$ws1.Range("C90").Value = 1.25
$ws1.Range("K90").Value = "10/25,27/2023"

$ws1.Range("A91").Value = 45231
$ws1.Range("B91").Value = "SL(2-0-0)"
$ws1.Range("H91").Value = 2
$ws1.Range("K91").Value = "11/10,30/2023"

$ws1.Range("A92").Value = 45261

# Row 93 becomes a "2024" year-marker row (style copied from an existing
# year-marker row) with a long blank placeholder remark in column F.
$ws1.Range("A10").Copy()
$ws1.Range("A93").PasteSpecial(-4122)
$ws1.Range("A93").Value = "2024"

$ws1.Range("D9").Copy()
$ws1.Range("B93").PasteSpecial(-4122)
$ws1.Range("D9").Copy()
$ws1.Range("D93").PasteSpecial(-4122)
$ws1.Range("D9").Copy()
$ws1.Range("F93").PasteSpecial(-4122)
$ws1.Range("F93").Value = "                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                                       "

$ws1.Range("A94").Value = 45292
$ws1.Range("A95").Value = 45323
$ws1.Range("A96").Value = 45352
$ws1.Range("A97").Value = 45383
$ws1.Range("A98").Value = 45413
$ws1.Range("A99").Value = 45444
$ws1.Range("A100").Value = 45474
$ws1.Range("A101").Value = 45505

# Extend Table13 to cover the new last row (103) before touching 102/103
$lo = $ws1.ListObjects.Item(1)
$lo.Resize($ws1.Range("A8:K103"))

# Preserve the old (bottom-border) formatting of row 102 by moving it to the
# new row 103 first ...
$ws1.Range("A102:K102").Copy($ws1.Range("A103:K103"))
$ws1.Range("G103").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'
$ws1.Range("A103").Value = 45566

# ... then give row 102 the regular (non-bottom-border) row formatting that
# row 101 has, and fill in its date.
$ws1.Range("A101:K101").Copy($ws1.Range("A102:K102"))
$ws1.Range("G102").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'
$ws1.Range("A102").Value = 45536

# ---------------------------------------------------------------------------
# Active tab: "2018 LEAVE CREDITS" becomes the selected / active sheet again
# ---------------------------------------------------------------------------
$ws1.Select()
$ws1.Range("F94").Select()

Write-Host "edits applied"
